$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The task list moved the "On-screen debug text" row from the bottom of the
# Engine group (row 11) to the top (row 2); the other rows in the Engine /
# Rorn Pool block shift down to fill the gap; "Bounding volume occlusion"
# was renamed to "Bounding sphere occlusion".

$ws.Cells.Item(2, 1).Value = "Engine"
$ws.Cells.Item(2, 2).Value = "On-screen debug text"
$ws.Cells.Item(2, 3).Value = 21

$ws.Cells.Item(3, 1).Value = "Engine"
$ws.Cells.Item(3, 2).Value = "Bounding sphere occlusion"
$ws.Cells.Item(3, 3).Value = 4

$ws.Cells.Item(4, 1).Value = "Engine"
$ws.Cells.Item(4, 2).Value = "FSAA"
$ws.Cells.Item(4, 3).Value = 2

$ws.Cells.Item(5, 1).Value = "Rorn Pool"
$ws.Cells.Item(5, 2).Value = "Build basic app"
$ws.Cells.Item(5, 3).Value = 2

$ws.Cells.Item(6, 1).Value = "Rorn Pool"
$ws.Cells.Item(6, 2).Value = "Game initialisation (position of balls, etc)"
$ws.Cells.Item(6, 3).Value = 14

$ws.Cells.Item(7, 1).Value = "Engine"
$ws.Cells.Item(7, 2).Value = "Basics of a physics engine"
$ws.Cells.Item(7, 3).Value = 21

$ws.Cells.Item(8, 1).Value = "Rorn Pool"
$ws.Cells.Item(8, 2).Value = "Player can move the cue ball (when appropriate)"
$ws.Cells.Item(8, 3).Value = 5

$ws.Cells.Item(9, 1).Value = "Rorn Pool"
$ws.Cells.Item(9, 2).Value = "Player can move cue"
$ws.Cells.Item(9, 3).Value = 5

$ws.Cells.Item(10, 1).Value = "Rorn Pool"
$ws.Cells.Item(10, 2).Value = "Player can take a shot"
$ws.Cells.Item(10, 3).Value = 10

$ws.Cells.Item(11, 1).Value = "Engine"
$ws.Cells.Item(11, 2).Value = "Compiled shaders"
$ws.Cells.Item(11, 3).Value = 5

# The review comment on "Compiled shaders" follows that task to its new
# row (B11).
$oldComment = $ws.Cells.Item(10, 2).Comment
$commentText = $oldComment.Text()
$oldComment.Delete()
$ws.Cells.Item(11, 2).AddComment($commentText)

# Update the active selection to match the new authored state.
$ws.Range("D2").Select()
